# Reorders the observation rows 2..10 in the active worksheet (rows are
# permuted as a whole, each row's full record moving together), and
# refreshes the AC/AJ/AK/AO "comment / substrate" columns to match their
# new row position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New content for each row, keyed by destination row number.
# (Each hashtable is the *entire* set of columns that the diff touches for
# that row; columns not listed here are left untouched because they were
# already identical across all rows.)
$rows = @{
    2  = @{ A = 80448769; B = 77506;  C = "Ovaliderad";                                   E = 6425;   F = "Garnlav";           G = "Alectoria sarmentosa";    H = "(Ach.) Ach.";      Q = 422991.0759451608; R = 6752021.173145968; AC = "Rikligt, hkb" }
    3  = @{ A = 80448775; B = 77506;  C = "Ovaliderad";                                   E = 6425;   F = "Garnlav";           G = "Alectoria sarmentosa";    H = "(Ach.) Ach.";      Q = 423036.1594514723; R = 6752009.000504656; AC = "Rikligt, hkb" }
    4  = @{ A = 80448777; B = 77506;  C = "Ovaliderad";                                   E = 6425;   F = "Garnlav";           G = "Alectoria sarmentosa";    H = "(Ach.) Ach.";      Q = 423115.1561234437; R = 6752009.239606674; AC = $null }
    5  = @{ A = 80448771; B = 73693;  C = "Ovaliderad";                                   E = 6440;   F = "Vitgrynig nållav";  G = "Chaenotheca subroscida";  H = "(Eitner) Zahlbr."; Q = 423289.9356373397; R = 6752041.978126496; AC = $null }
    6  = @{ A = 80448779; B = 77506;  C = "Ovaliderad";                                   E = 6425;   F = "Garnlav";           G = "Alectoria sarmentosa";    H = "(Ach.) Ach.";      Q = 422962.8083476268; R = 6752021.785183201; AC = "Rikligt" }
    7  = @{ A = 80448780; B = 77506;  C = "Ovaliderad";                                   E = 6425;   F = "Garnlav";           G = "Alectoria sarmentosa";    H = "(Ach.) Ach.";      Q = 423056.1482692101; R = 6751963.779848268; AC = "Spritt" }
    8  = @{ A = 80448772; B = 81236;  C = "Ovaliderad";                                   E = 1312;   F = "Gammelgransskål";   G = "Pseudographis pinicola";  H = "(Nyl.) Rehm";      Q = 423289.9356373397; R = 6752041.978126496; AC = $null }
    9  = @{ A = 80448773; B = 77506;  C = "Ovaliderad";                                   E = 6425;   F = "Garnlav";           G = "Alectoria sarmentosa";    H = "(Ach.) Ach.";      Q = 422635.9957601223; R = 6751949.037152009; AC = $null; AJ = "vanlig tall"; AK = "Pinus sylvestris var. sylvestris"; AO = "Pinus sylvestris var. sylvestris" }
    10 = @{ A = 80448778; B = 56395;  C = "Godkänd baserat på observatörens uppgifter";   E = 100109; F = "Tretåig hackspett"; G = "Picoides tridactylus";    H = "(Linnaeus, 1758)"; Q = 423115.1561234437; R = 6752009.239606674; AC = $null; AJ = $null; AK = $null; AO = $null }
}

foreach ($r in $rows.Keys) {
    $data = $rows[$r]

    $ws.Range("A$r").Value = $data.A
    $ws.Range("B$r").Value = $data.B
    $ws.Range("C$r").Value = $data.C
    $ws.Range("E$r").Value = $data.E
    $ws.Range("F$r").Value = $data.F
    $ws.Range("G$r").Value = $data.G
    $ws.Range("H$r").Value = $data.H
    $ws.Range("Q$r").Value = $data.Q
    $ws.Range("R$r").Value = $data.R

    if ($data.AC -eq $null) {
        $ws.Range("AC$r").ClearContents()
    } else {
        $ws.Range("AC$r").Value = $data.AC
    }

    if ($data.ContainsKey("AJ")) {
        if ($data.AJ -eq $null) { $ws.Range("AJ$r").ClearContents() } else { $ws.Range("AJ$r").Value = $data.AJ }
    }
    if ($data.ContainsKey("AK")) {
        if ($data.AK -eq $null) { $ws.Range("AK$r").ClearContents() } else { $ws.Range("AK$r").Value = $data.AK }
    }
    if ($data.ContainsKey("AO")) {
        if ($data.AO -eq $null) { $ws.Range("AO$r").ClearContents() } else { $ws.Range("AO$r").Value = $data.AO }
    }
}
